$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# The old EPA "MortalityRiskValuation" URL moved; point the cell + hyperlink
# at the new https://www.epa.gov/... address (same target used for the
# "whatvalue" in-page anchor).
$newUrl = "https://www.epa.gov/environmental-economics/mortality-risk-valuation#whatvalue"
$newBase = "https://www.epa.gov/environmental-economics/mortality-risk-valuation"

$cell = $ws.Range("B6")
$cell.Value = $newUrl

# Turn it into a real clickable hyperlink (Address + in-page anchor), the
# way the previous URL-only text cell was not wired up.
$ws.Hyperlinks.Add($cell, $newBase, "whatvalue", "", "$newBase - whatvalue")

# Hyperlinks.Add re-applies the built-in Hyperlink cell style (adding a new
# duplicate style record along the way); restore the existing Hyperlink
# style so the cell keeps using the workbook's original style entry.
$cell.Style = "Hyperlink"

# Enable iterative calculation (used by the model elsewhere in the workbook).
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.00001
